# Ressourcenplanung Online Shop Projekt
# Applies the resource-planning update to the "Netzplan" (network plan) sheet:
#  - work packages in rows 6, 9, 10, 11, 14, 16 are now staffed by 2 people
#    instead of 1 (column H, "Personen")
#  - the overall staffing-efficiency factor in I12 drops from 100% to 45%
#  - a handful of previously hard-coded helper cells in the Gantt/early-late
#    time grid become real formulas so the sheet recalculates consistently
#  - two formulas are tidied up (AZ4, X17)
#  - view settings (zoom / top-left cell / selection) are refreshed on both
#    sheets, and the workbook's tab-bar ratio changes

$wb = $excel.ActiveWorkbook

$psp = $wb.Worksheets.Item("PSP")
$np  = $wb.Worksheets.Item("Netzplan")

# ---------------------------------------------------------------------------
# 1. Root inputs: staffing counts and efficiency factor
# ---------------------------------------------------------------------------
$np.Range("H6").Value  = "2"
$np.Range("H9").Value  = "2"
$np.Range("H10").Value = "2"
$np.Range("H11").Value = "2"
$np.Range("H14").Value = "2"
$np.Range("H16").Value = "2"

$np.Range("I12").Value = 0.45

# ---------------------------------------------------------------------------
# 2. Per-person capacity-adjustment formulas (only the two-person packages
#    get the coordination-overhead correction)
# ---------------------------------------------------------------------------
$np.Range("I6").Formula  = "=(100%+(100%-`$I`$12))/H6"
$np.Range("I9").Formula  = "=(100%+(100%-`$I`$12))/H9"
$np.Range("I10").Formula = "=(100%+(100%-`$I`$12))/H10"
$np.Range("I11").Formula = "=(100%+(100%-`$I`$12))/H11"
$np.Range("I14").Formula = "=(100%+(100%-`$I`$12))/H14"
$np.Range("I16").Formula = "=(100%+(100%-`$I`$12))/H16"

# ---------------------------------------------------------------------------
# 3. Early-time grid (row 2 block) - turn hard-coded numbers into formulas
# ---------------------------------------------------------------------------
$np.Range("N2").Formula  = "=L2+L4"
$np.Range("S2").Formula  = "=Q2+Q4"
$np.Range("V2").Formula  = "=S2"
$np.Range("X2").Formula  = "=V2+V4"
$np.Range("AF2").Formula = "=MAX(X2,AC8,X14)"
$np.Range("AH2").Formula = "=AF2+AF4"

# ---------------------------------------------------------------------------
# 4. Late-time grid (row 5 block)
# ---------------------------------------------------------------------------
$np.Range("X5").Formula  = "=AF5"
$np.Range("BL5").Formula = "=BO5"

# ---------------------------------------------------------------------------
# 5. Work-package block starting row 8 (AC8 branch)
# ---------------------------------------------------------------------------
$np.Range("X8").Formula  = "=V8+V10"
$np.Range("AA8").Formula = "=X8"
$np.Range("AC8").Formula = "=AA8+AA10"

# ---------------------------------------------------------------------------
# 6. Row 11 branch
# ---------------------------------------------------------------------------
$np.Range("X11").Formula  = "=AA11"
$np.Range("AC11").Formula = "=AF5"

# ---------------------------------------------------------------------------
# 7. Work-package block starting row 14 (AC14 branch)
# ---------------------------------------------------------------------------
$np.Range("Q14").Formula  = "=N2"
$np.Range("S14").Formula  = "=Q14+Q16"
$np.Range("V14").Formula  = "=S14"
$np.Range("X14").Formula  = "=V14+V16"
$np.Range("AA14").Formula = "=X14"
$np.Range("AC14").Formula = "=AA14+AA16"

# ---------------------------------------------------------------------------
# 8. Row 17 branch
# ---------------------------------------------------------------------------
$np.Range("AC17").Formula = "=BO5"

# ---------------------------------------------------------------------------
# 9. Formula tidy-ups
# ---------------------------------------------------------------------------
$np.Range("AZ4").Formula = "=J14"
$np.Range("X17").Formula = "=MIN(AA17,AF5)"

# ---------------------------------------------------------------------------
# 10. View settings
# ---------------------------------------------------------------------------
$psp.Range("A1").Select()
$psp.Application.ActiveWindow.Zoom = 100
$psp.Application.ActiveWindow.ScrollRow = 1
$psp.Application.ActiveWindow.ScrollColumn = 1

$np.Range("A1").Select()
$np.Application.ActiveWindow.Zoom = 100
$np.Application.ActiveWindow.ScrollRow = 10
$np.Application.ActiveWindow.ScrollColumn = 1

$wb.Application.ActiveWindow.TabRatio = 60

$wb.Save()
